# fix sheet name of designation and department xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Employees" to "Designation"
$ws.Name = "Designation"

# Add the new row with the "Project Manager" value
$ws.Range("A2").Value = "Project Manager"

# Update selection to match the target state
$ws.Range("C9").Select()
